$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header cells: "..._old" -> "..._FV2304", "..._new" -> "..._FV2310"
# (used in the AHB-Diff header row, A1:U1)
for ($col = 1; $col -le 21; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $val = $cell.Value2
    if ($val -ne $null) {
        $newVal = $val -replace '_old$', '_FV2304'
        $newVal = $newVal -replace '_new$', '_FV2310'
        if ($newVal -ne $val) {
            $cell.Value = $newVal
        }
    }
}

# Turn the data range into an Excel Table ("Table1") with autofilter
$dataRange = $ws.Range("A1:U63")
$table = $ws.ListObjects.Add(1, $dataRange, $null, 1)
$table.Name = "Table1"
$table.TableStyle = ""

# Freeze the header row (row 1)
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
